# Doing Updates for Financials
# Insert a new "most recent period" column before column D on the MCBC sheet,
# shifting the existing D:K data right to E:L, and populate the new column D
# with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCBC")

# Insert a new blank column at D; this shifts D:K -> E:L automatically
# (dimension, row spans, etc. are all adjusted by Excel).
$ws.Columns("D:D").Insert()

# Give the new column D the same formatting as column E (which holds what used
# to be column D), so number formats / styles line up with the rest of the data.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate new column D with the latest period's values ---

# Income Statement
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 69000
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 9900
$ws.Cells.Item(18, 4).Value = 59200
$ws.Cells.Item(20, 4).Value = -26800
$ws.Cells.Item(21, 4).Value = 34900
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 4).Value = 32400
$ws.Cells.Item(24, 4).Value = 6000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 26400
$ws.Cells.Item(27, 4).Value = 26400
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 26800
$ws.Cells.Item(33, 4).Value = 26400
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 26400

# Balance Sheet
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 40500
$ws.Cells.Item(42, 4).Value = 142300
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 44900
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 3400
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 1975100
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(59, 4).Value = 6300
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(61, 4).Value = 41200
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 1784300
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = -24700
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 190900
$ws.Cells.Item(77, 4).Value = 0

# Cash Flow Statement
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 26400
$ws.Cells.Item(83, 4).Value = 2600
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 29600
$ws.Cells.Item(91, 4).Value = -1300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -77200
$ws.Cells.Item(96, 4).Value = -8500
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 57400
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = 9800

# Row 49 (Goodwill) previously showed "NA" across every period; the source
# data now reports 0 for all of the historical periods too (columns D-J),
# while the oldest period (now column K) keeps its value of 100.
for ($c = 4; $c -le 10; $c++) {
    $ws.Cells.Item(49, $c).Value = 0
}
